$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("V3").Value = 1.58
$ws.Range("R4").Value = 1.72
$ws.Range("R5").Value = 1.54
$ws.Range("M6").Value = 10.9
$ws.Range("G7").Value = 1.36
$ws.Range("H7").Value = 4.65
$ws.Range("I7").Value = 7.2
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.52
$ws.Range("R7").Value = 2.22
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 1.85
$ws.Range("W7").Value = 8.25
$ws.Range("X7").Value = 7.2
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 10.5
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 15
$ws.Range("AD7").Value = 9.5
$ws.Range("AE7").Value = 18
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 23
$ws.Range("AL7").Value = 60
$ws.Range("AM7").Value = 500
$ws.Range("AN7").Value = 3.3
$ws.Range("AP7").Value = 14.5
$ws.Range("AU7").Value = 7.8
$ws.Range("AY7").Value = 35
$ws.Range("BB7").Value = 400
$ws.Range("Q11").Value = 1.92
$ws.Range("R11").Value = 1.82
$ws.Range("AT12").Value = 2.62
$ws.Range("AT13").Value = 2.62
$ws.Range("AT14").Value = 2.62
$ws.Range("AT21").Value = 2.37

